$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 36, continuing the data table (same layout/style as row 35)
$srcRow = $ws.Range("A35:E35")
$dstRow = $ws.Range("A36:E36")

# Copy formatting (and formats) from the last existing row down to the new row
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the values for the new record (29 May 2018)
$ws.Range("A36").Value = Get-Date -Year 2018 -Month 5 -Day 29 -Hour 0 -Minute 0 -Second 0
$ws.Range("B36").Value = "Di Giacomo Caterina"
$ws.Range("C36").Value = "Tela Leggera"
$ws.Range("D36").Value = "Mt."
$ws.Range("E36").Value = 26
